# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp text (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 22:22"

# --- Refresh numeric data for several countries ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 557043
$ws.Range("C4").Value = 24164
$ws.Range("E4").Value = 503722
$ws.Range("F4").Value = 11766
$ws.Range("G4").Value = 1375
$ws.Range("H4").Value = 21952

# Row 5: Espana
$ws.Range("B5").Value = 166127
$ws.Range("C5").Value = 3100
$ws.Range("E5").Value = 86623
$ws.Range("G5").Value = 507
$ws.Range("H5").Value = 17113

# Row 8: Alemania
$ws.Range("B8").Value = 127459
$ws.Range("C8").Value = 2007
$ws.Range("E8").Value = 64163
$ws.Range("G8").Value = 125
$ws.Range("H8").Value = 2996

# Row 17: Brasil
$ws.Range("B17").Value = 21066
$ws.Range("C17").Value = 104
$ws.Range("E17").Value = 19747
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 1146

# Row 31: Noruega
$ws.Range("E31").Value = 6325
$ws.Range("G31").Value = 9
$ws.Range("H31").Value = 128

# --- Reorder "Birmania" / "Puerto Rico" (swap their table position) and
#     refresh the data that moved with them ---
# Row 154 becomes Birmania (with fresh data); Row 155 becomes Puerto Rico
# (keeping the data previously shown for Puerto Rico on row 154).
$ws.Range("A154").Value = "Birmania"
$ws.Range("B154").Value = 41
$ws.Range("C154").Value = 3
$ws.Range("D154").Value = 2
$ws.Range("E154").Value = 35
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 4

$ws.Range("A155").Value = "Puerto Rico"
$ws.Range("B155").Value = 39
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 1
$ws.Range("E155").Value = 36
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 2
